# Updated cryptos list on Sat Nov  4 04:39:01 UTC 2023 with GitHub Actions
#
# The sheet is a "Coin / Link / Price / Volume(1h)" table (columns B-E,
# rows 2-51). Every data cell is stored as text (even the ones that look
# like plain numbers, e.g. "41.92" or "1.07" - those are prices such as
# "34.973.55" with thousands separators, so the whole Price column is
# text). Plain Excel .Value assignment auto-detects numeric-looking
# strings and stores them as numbers, which would change the cell type,
# so Set-TextCell forces the cell to Text first, assigns the value, then
# drops back to the default "Normal" style so no stray NumberFormat is
# left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($sheet, $addr, $val) {
    $cell = $sheet.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# --- Price / Volume(1h) refresh ------------------------------------------
$ws.Range("D2").Value = "34.954.40"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "1.844.72"
$ws.Range("E3").Value = "  +2.03%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  +0.98%  "
$ws.Range("E6").Value = "  +2.49%  "
Set-TextCell $ws "D8" "41.69"
$ws.Range("E8").Value = "  +7.00%  "
Set-TextCell $ws "D9" "0.328"
$ws.Range("E9").Value = "  +3.18%  "
Set-TextCell $ws "D10" "0.0693"
$ws.Range("E10").Value = "  +2.19%  "
$ws.Range("E11").Value = "  -0.80%  "
$ws.Range("D12").Value = "2.113.52"
$ws.Range("E12").Value = "  +2.11%  "
Set-TextCell $ws "D13" "11.43"
$ws.Range("E13").Value = "  +5.43%  "
$ws.Range("D14").Value = "1.850.93"
$ws.Range("E14").Value = "  +2.18%  "
Set-TextCell $ws "D15" "0.673"
$ws.Range("E15").Value = "  +2.13%  "
$ws.Range("E16").Value = "  +3.04%  "
$ws.Range("D17").Value = "34.966.18"
$ws.Range("E17").Value = "  +0.34%  "
Set-TextCell $ws "D18" "69.98"
$ws.Range("E18").Value = "  +1.14%  "
$ws.Range("D19").Value = "0.0₃0791"
$ws.Range("E19").Value = "  +1.53%  "
Set-TextCell $ws "D20" "240.61"
$ws.Range("E20").Value = "  +0.86%  "
Set-TextCell $ws "D21" "12.21"
$ws.Range("E21").Value = "  +4.21%  "
Set-TextCell $ws "D22" "4.75"
$ws.Range("E22").Value = "  +2.69%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("E24").Value = "  +1.11%  "
Set-TextCell $ws "D25" "172.31"
$ws.Range("E25").Value = "  -0.71%  "
Set-TextCell $ws "D26" "7.80"
$ws.Range("E26").Value = "  +0.77%  "
Set-TextCell $ws "D27" "17.50"
$ws.Range("E27").Value = "  +2.00%  "
$ws.Range("E28").Value = "  +4.21%  "
Set-TextCell $ws "D29" "1.65"
$ws.Range("E29").Value = "  +9.59%  "
$ws.Range("E30").Value = "  +0.06%  "
Set-TextCell $ws "D31" "0.0553"
$ws.Range("E31").Value = "  +1.51%  "
$ws.Range("E32").Value = "  -0.32%  "
Set-TextCell $ws "D33" "3.91"
$ws.Range("E33").Value = "  -0.12%  "
$ws.Range("E34").Value = "  +21.97%  "
$ws.Range("E35").Value = "  +11.79%  "
Set-TextCell $ws "D36" "0.746"
$ws.Range("E36").Value = "  +9.37%  "
$ws.Range("E37").Value = "  +4.77%  "
Set-TextCell $ws "D38" "1.06"
$ws.Range("E38").Value = "  +11.69%  "
Set-TextCell $ws "D39" "89.72"
$ws.Range("D40").Value = "1.349.85"
$ws.Range("E40").Value = "  +3.23%  "
$ws.Range("E41").Value = "  +2.78%  "
Set-TextCell $ws "D42" "14.60"
$ws.Range("E42").Value = "  +3.54%  "
Set-TextCell $ws "D43" "2.28"
$ws.Range("E43").Value = "  +4.60%  "
Set-TextCell $ws "D45" "2.75"
$ws.Range("E45").Value = "  +2.16%  "
$ws.Range("E46").Value = "  +4.08%  "
$ws.Range("E47").Value = "  +3.96%  "
$ws.Range("D48").Value = "2.033.96"
$ws.Range("E48").Value = "  +1.80%  "
$ws.Range("E51").Value = "  -0.52%  "

# --- Rows 49/50 swap: THORChain overtakes PaxDollar in the ranking -------
$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextCell $ws "D49" "3.41"
$ws.Range("E49").Value = "  +15.92%  "

$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextCell $ws "D50" "1.01"
$ws.Range("E50").Value = "  +0.11%  "
